# Updates cryptos list values (Price / Volume(1h)) for rows 2-51
# generated from the commit diff - GitHub Actions cryptos list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.315.90'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '2.245.93'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''306.61'
$ws.Range('E5').Value = '  -2.41%  '
$ws.Range('D6').Value = '''96.06'
$ws.Range('E6').Value = '  -2.55%  '
$ws.Range('D7').Value = '''0.573'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('D10').Value = '''34.97'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').Value = '''0.0813'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = '''7.26'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').Value = '2.341.55'
$ws.Range('E14').Value = '  +5.02%  '
$ws.Range('D15').Value = '2.588.36'
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').Value = '''0.834'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '''13.59'
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('D18').Value = '44.117.04'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = '0.0₃0967'
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('D20').Value = '''6.38'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').Value = '''12.12'
$ws.Range('E21').Value = '  -7.09%  '
$ws.Range('D22').Value = '''65.65'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').Value = '''238.37'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('D24').Value = '''2.95'
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '''38.66'
$ws.Range('E27').Value = '  +6.02%  '
$ws.Range('D28').Value = '''9.95'
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('E29').Value = '  +2.79%  '
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('D31').Value = '''5.90'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = '''152.70'
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('E36').Value = '  +2.16%  '
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('D38').Value = '''1.76'
$ws.Range('E38').Value = '  -7.33%  '
$ws.Range('D39').Value = '''3.57'
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('D40').Value = '''3.86'
$ws.Range('E40').Value = '  -3.32%  '
$ws.Range('D41').Value = '''14.48'
$ws.Range('E41').Value = '  -6.68%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').Value = '1.754.16'
$ws.Range('E44').Value = '  +2.80%  '
$ws.Range('D45').Value = '''82.98'
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('D46').Value = '''0.191'
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').Value = '''100.21'
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').Value = '''4.96'
$ws.Range('E48').Value = '  -3.24%  '
$ws.Range('D49').Value = '''8.12'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').Value = '''54.82'
$ws.Range('E51').Value = '  -2.10%  '
